# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# D3: Good Roaming Calculation (%) updated from 98.2 to 98
$ws.Range("D3").Value = 98

# E12: Driver Vintage date populated (was blank)
# Write as a literal text value (not an auto-converted date serial), matching
# the plain-text date strings already used elsewhere in column E (e.g. E13/E14).
$ws.Range("E12").Formula = '="2022-08-29"'
$ws.Range("E12").Copy($null)
$ws.Range("E12").PasteSpecial(-4163)  # xlPasteValues: bake the formula result into a static value
